$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vessels")

# Weather restrictions constraint: update transfer_time (F) and max_time_offshore (H)
$ws.Range("F2").Value = 0.25
$ws.Range("H2").Value = 12

$ws.Range("F3").Formula = "=20/60"
$ws.Range("H3").Value = 12

$ws.Range("F4").Formula = "=30/60"
$ws.Range("H4").Value = 24

# Update the active selection to match the saved view state
$ws.Range("H7").Select()
